# "tipos de experimentos correcto"
# Fix the Google Drive file id in the "Tipos de Experimentos" link stored
# in cell D7 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("D7")

$oldId = "10A4EJFGWfebygOzYqxJxFSMew-zIN9r7"
$newId = "10NvXg9FUUU0muewzUiHE0CgGXWOo_wpK"

[string]$text = $cell.Value()
$text = $text.Replace($oldId, $newId)
$cell.Value = $text
